# "modif plan test": renumber rows 6-7 (insert of a new item shifted the
# existing #4/#5 down to #5/#6) and update the saved view/selection so the
# workbook reopens scrolled to, and focused on, the newly numbered row 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

$ws.Activate()

# Scroll the window so row 6 becomes the top visible row (sheetView
# topLeftCell="A6") and leave the selection on A7.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A7").Select()
